# Auto update stock data
# Updates the "Date_1" column (A) to the new refresh date for every data
# row, and refreshes the "EBITDA" column (B) values that moved, while
# keeping everything text (these columns store text-formatted numbers /
# dates, not real numeric/date cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new EBITDA (B column) value; $null means "date only, no B change"
$rows = [ordered]@{
    2  = "5.45"
    8  = "7.83"
    14 = "2.89"
    20 = "12.45"
    26 = "10.35"
    32 = "26.21"
    38 = $null
    44 = "11.46"
    50 = "12.16"
    56 = "35.69"
    62 = "12.20"
    68 = "13.98"
    74 = "16.78"
}

foreach ($r in $rows.Keys) {
    # Leading apostrophe keeps the written value as literal text (matches
    # the existing text-stored dates/numbers instead of Excel coercing
    # "2025/12/09" into a real date serial or "5.45" into a float).
    $ws.Cells.Item($r, 1).Value = "'2025/12/09"

    $newB = $rows[$r]
    if ($null -ne $newB) {
        $ws.Cells.Item($r, 2).Value = "'" + $newB
    }
}
